# Regenerate the localization status report for archive:
#  - flip the "Ready for handoff" status to "In Translation" everywhere it
#    appears (Overview summary columns + the per-locale Status column)
#  - re-fit the width of the columns that held the old, longer status text

$wb = $excel.ActiveWorkbook

# --- Update status text -----------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Re-fit the affected Status columns now that the text is shorter ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
